# Fill in the previously-blank "role" indicator cells (projectManager,
# frontendDeveloper, backendDeveloper, dataScientist, dataEngineer =
# columns G:K) with explicit 0 values for the rows where they were missing,
# correct the predecessor for row 5 / row 6, and leave the active selection
# on H5 (matching the author's final click position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (task A)
$ws.Range("H2:K2").Value = 0

# Row 3 (task B)
$ws.Range("H3:K3").Value = 0

# Row 4 (task C)
$ws.Range("H4:K4").Value = 0

# Row 5 (task D) - add missing predecessor + zero-fill role columns
$ws.Range("C5").Value = "A"
$ws.Range("E5").Value = 0
$ws.Range("G5:K5").Value = 0

# Row 6 (task D1) - predecessor corrected from "A" to "D1"
$ws.Range("C6").Value = "D1"

# Row 8 (task D3)
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0

# Row 9 (task D4.0)
$ws.Range("G9:K9").Value = 0

# Row 10 (task D4.1)
$ws.Range("G10").Value = 0
$ws.Range("I10:K10").Value = 0

# Row 11 (task D4.2)
$ws.Range("G11:H11").Value = 0
$ws.Range("J11:K11").Value = 0

# Row 12 (task D4.3)
$ws.Range("G12:I12").Value = 0
$ws.Range("K12").Value = 0

# Row 13 (task D4.4)
$ws.Range("G13:J13").Value = 0

# Row 14 (task D4.5)
$ws.Range("G14:K14").Value = 0

# Row 16 (task D6)
$ws.Range("G16").Value = 0
$ws.Range("J16").Value = 0

# Row 17 (task D7)
$ws.Range("J17").Value = 0

# Row 19 (task E)
$ws.Range("H19:K19").Value = 0

# Row 20 (task F)
$ws.Range("H20:K20").Value = 0

# Row 21 (task G)
$ws.Range("H21:K21").Value = 0

# Row 22 (task H)
$ws.Range("H22:K22").Value = 0

# Restore the active cell/selection to H5, as left by the editor.
$ws.Range("H5").Select()
